$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "4787 ms"
$ws.Range("E3").Value = "4468 ms"
$ws.Range("E4").Value = "4418 ms"
$ws.Range("E5").Value = "4294 ms"
$ws.Range("E6").Value = "4273 ms"
$ws.Range("E7").Value = "4235 ms"
